$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the run(s) covering [start,end) away from their neighbours by
# briefly toggling a character property on then back off. Word (and this
# interop) splits a run whenever a sub-range gets a differing formatting
# value; flipping Bold on then off again leaves the *visible* formatting
# untouched (Bold=False is the already-inherited default, so the temporary
# <w:b/> is dropped again) while still forcing the run boundary to land
# exactly at start/end.
# ---------------------------------------------------------------------------
function Split-RunBoundary($doc, $start, $end) {
    $r = $doc.Range($start, $end)
    $r.Bold = 1
    $r.Bold = 0
}

# ===========================================================================
# Change 1: "the grid of strings" paragraph - explain the array sizing fix.
# ===========================================================================
$oldText1 = " read a CSV file and translate it into a grid of strings. The first possible limitation of my project is the size of the grid of strings. For the sake of the example, I chose to set the size as 3 by 3, however this would mean there would be issues if the inputted CSV file was larger than the example one. I did also contemplate doing a 10 by 10, but that would sacrifice performance, as the for loops would loop through many null elements when values aren't present."

$newPart1a = " read a CSV file and translate it into a grid of strings. The first possible limitation of my project is the size of the grid of strings. "
$newPart1b = "I amended the code so that it would first parse the file to obtain size of the grid to create an array[][] of the correct size, but the limitation is that this required the file to be parsed twice, which could affect performance. However, I felt in the long run, it may be better to parse the file twice, than to use the other method of copying an array to a larger array for every iteration after the initial size."

$newText1 = $newPart1a + $newPart1b

$d.Content.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 2) | Out-Null

$t = $d.Content.Text
$start1 = $t.IndexOf($newPart1a)
$splitAt1 = $start1 + $newPart1a.Length
$end1 = $splitAt1 + $newPart1b.Length
Split-RunBoundary $d $splitAt1 $end1

# ===========================================================================
# Change 2: "calculateEquations" paragraph - floats -> doubles, trim operator
# clause duplication.
# ===========================================================================
$oldText2 = " method, which first ensures a string isn't empty, null or just a blank space and then splits it into an array of smaller strings, before choosing between whether it's Cell information, in which it will replace the cell with the cells value and place that in a stack of floats, whether it's an operator, in which it will perform a calculation and place the value within the stack of floats, or whether it's a number, in which it will parse the float and place It in the stack of floats. At the end there should be only one value, which will be formatted to make sure it doesn't show a decimal unless necessary."

$newPart2a = " method, which first ensures a string isn't empty, null or just a blank space and then splits it into an array of smaller strings, before choosing between whether it's Cell information, in which it will replace the cell with the cells value and place that in a stack of "
$newPart2b = "doubles"
$newPart2c = ", or whether it's a number, in which it will parse the "
$newPart2d = "double"
$newPart2e = " and place It in the stack of "
$newPart2f = "doubles"
$newPart2g = ". At the end there should be only one value, which will be formatted to make sure it doesn't show a decimal unless necessary."

$newText2 = $newPart2a + $newPart2b + $newPart2c + $newPart2d + $newPart2e + $newPart2f + $newPart2g

$d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2) | Out-Null

$t = $d.Content.Text
$start2 = $t.IndexOf($newPart2a)

$p0 = $start2
$p1 = $p0 + $newPart2a.Length
$p2 = $p1 + $newPart2b.Length
$p3 = $p2 + $newPart2c.Length
$p4 = $p3 + $newPart2d.Length
$p5 = $p4 + $newPart2e.Length
$p6 = $p5 + $newPart2f.Length
$p7 = $p6 + $newPart2g.Length

Split-RunBoundary $d $p1 $p2
Split-RunBoundary $d $p2 $p3
Split-RunBoundary $d $p3 $p4
Split-RunBoundary $d $p4 $p5
Split-RunBoundary $d $p5 $p6
Split-RunBoundary $d $p6 $p7
